$wb = $excel.ActiveWorkbook

# --- zh-cn sheet (row 6 corresponds to 841026d0-8eff-4a88-88d8-a39b1680ae53) ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to fit the new long error message.
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# Fill in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# and "Error Detail" columns for the handback report that was just generated.
$wsZh.Range("I6").Value = "841026d0-8eff-4a88-88d8-a39b1680ae53.md"
$wsZh.Range("J6").Value = "841026d0-8eff-4a88-88d8-a39b1680ae53.28b072df7ae65089e701c4c793600a7ac536773c.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-08-26 22:42:08"
$wsZh.Range("P6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37b08580129ec248456f99f0ad3de84a740688bc/e2e/841026d0-8eff-4a88-88d8-a39b1680ae53.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c720f50dfca02b061fe6c7414fac6631318500f/e2e/841026d0-8eff-4a88-88d8-a39b1680ae53.md."

$wsZh.Hyperlinks.Add($wsZh.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c720f50dfca02b061fe6c7414fac6631318500f/e2e/841026d0-8eff-4a88-88d8-a39b1680ae53.md", "", "", "841026d0-8eff-4a88-88d8-a39b1680ae53.md")

# --- de-de sheet (row 6 corresponds to 841026d0-8eff-4a88-88d8-a39b1680ae53) ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$wsDe.Range("I6").Value = "841026d0-8eff-4a88-88d8-a39b1680ae53.md"
$wsDe.Range("J6").Value = "841026d0-8eff-4a88-88d8-a39b1680ae53.28b072df7ae65089e701c4c793600a7ac536773c.de-de.xlf"
$wsDe.Range("K6").Value = "2016-08-26 22:42:14"
$wsDe.Range("P6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37b08580129ec248456f99f0ad3de84a740688bc/e2e/841026d0-8eff-4a88-88d8-a39b1680ae53.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c720f50dfca02b061fe6c7414fac6631318500f/e2e/841026d0-8eff-4a88-88d8-a39b1680ae53.md."

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0c720f50dfca02b061fe6c7414fac6631318500f/e2e/841026d0-8eff-4a88-88d8-a39b1680ae53.md", "", "", "841026d0-8eff-4a88-88d8-a39b1680ae53.md")
